# Apply the "First set of budget test working." edit:
#  - Correct the Programme/Cost-Centre style data in row 2 (A2, B2)
#  - Fix the stray full-row selection left on the sheet (A4:XFD4 -> A4)
#  - Give column B an explicit best-fit width, matching the narrower
#    numeric content now shown there

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budgets")

# --- Row 2 data corrections ---------------------------------------------
$ws.Range("A2").Value = 109076
$ws.Range("B2").Value = 11272001

# --- Column B best-fit width --------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 8.33

# --- Selection fix: only A4 should be selected, not the whole row -------
$ws.Range("A4").Select()

$wb.Save()
